$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "STM32L152RE"
$ws.Range("G5").Value = 0.010296
$ws.Range("H5").Value = 1.0361910000000001
$ws.Range("I5").Formula = '=($D$16/H5)/1000000000'
$ws.Range("J5").Formula = '=I5/G5'

$ws.Range("H24").Select()
